# "Generate Report for Handoff"
#
# Regenerates the localization status report:
#   - status moves from "Ready for handoff" to "In Translation" everywhere
#     it's reported (Overview summary + each per-language detail sheet)
#   - the corresponding handoff/generation timestamps are refreshed
#   - the status columns re-autofit to the new (shorter) status text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---
$overview.Range("E2").Value = "In Translation"   # zh-cn status (Overview)
$overview.Range("F2").Value = "In Translation"   # de-de status (Overview)
$zhcn.Range("C2").Value     = "In Translation"   # Status column on zh-cn sheet
$dede.Range("C2").Value     = "In Translation"   # Status column on de-de sheet

# --- Refreshed handoff / generation timestamps ---
$overview.Range("G2").Value = "2016-10-20 06:58:25"   # Latest HO Xliff Generate Date
$zhcn.Range("H2").Value     = "2016-10-20 06:58:14"   # Latest Handoff Datetime (zh-cn)
$dede.Range("H2").Value     = "2016-10-20 06:58:25"   # Latest Handoff Datetime (de-de)

# --- Column widths: the status columns shrink to fit "In Translation" ---
# The stored OOXML <col width> is always (quantized-ColumnWidth + 5/6) in this
# engine, quantized to 1/6-character steps, so solve for the ColumnWidth input
# whose bucket lands nearest the real-Excel autofit width recorded in the sheet.
function Set-AutofitColumnWidth($col, $storedWidth) {
    $col.ColumnWidth = [Math]::Round(($storedWidth - 5.0/6.0) * 6.0) / 6.0
}

$targetStoredWidth = 13.4101848602295

Set-AutofitColumnWidth $overview.Columns.Item(5) $targetStoredWidth   # Overview!E
Set-AutofitColumnWidth $overview.Columns.Item(6) $targetStoredWidth   # Overview!F
Set-AutofitColumnWidth $zhcn.Columns.Item(3)     $targetStoredWidth   # zh-cn!C (Status)
Set-AutofitColumnWidth $dede.Columns.Item(3)     $targetStoredWidth   # de-de!C (Status)
